# Applies the "quarto template" touch-up to doc_template.docx:
#   * explicit Calibri rFonts on the Author/Date lines and three of the
#     FirstParagraph body lines (both the paragraph mark and the run)
#   * drop the stray "_GoBack" bookmark from the empty "section" heading
#     paragraph (Word renumbers the following "section-1" bookmark down
#     from 6 to 5 as a result)
#   * merge the two runs in the "Markdown code '#####'" paragraph into one

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $xml = "<w:p $wNs>$innerXml</w:p>"
    $r.InsertXML($xml)
}

$calibriRpr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>'

# --- Author: "Name" ---------------------------------------------------
$inner = '<w:pPr><w:pStyle w:val="Author"/>' + $calibriRpr + '</w:pPr>'
$inner = $inner + '<w:r>' + $calibriRpr + '<w:t>Name</w:t></w:r>'
Set-ParagraphXml 2 $inner

# --- Date: "8/16/2019" --------------------------------------------------
$inner = '<w:pPr><w:pStyle w:val="Date"/>' + $calibriRpr + '</w:pPr>'
$inner = $inner + '<w:r>' + $calibriRpr + '<w:t>8/16/2019</w:t></w:r>'
Set-ParagraphXml 3 $inner

# --- FirstParagraph: "Here is the text below a header." ----------------
$inner = '<w:pPr><w:pStyle w:val="FirstParagraph"/>' + $calibriRpr + '</w:pPr>'
$inner = $inner + '<w:r>' + $calibriRpr + '<w:t>Here is the text below a header.</w:t></w:r>'
Set-ParagraphXml 5 $inner

# --- FirstParagraph: "For our purposes, ..." ----------------------------
$inner = '<w:pPr><w:pStyle w:val="FirstParagraph"/>' + $calibriRpr + '</w:pPr>'
$inner = $inner + '<w:r>' + $calibriRpr + '<w:t>For our purposes, we will keep all headers the same (except #####).</w:t></w:r>'
Set-ParagraphXml 7 $inner

# --- FirstParagraph: "Additional text here" -----------------------------
$inner = '<w:pPr><w:pStyle w:val="FirstParagraph"/>' + $calibriRpr + '</w:pPr>'
$inner = $inner + '<w:r>' + $calibriRpr + '<w:t>Additional text here</w:t></w:r>'
Set-ParagraphXml 9 $inner

# --- Empty Heading5 paragraph holding the "section" bookmark: drop the
#     hidden "_GoBack" bookmark that Word had stamped at the last edit
#     position. The following bookmark ("section-1") is renumbered by
#     Word automatically once this one goes away.
$inner = '<w:pPr><w:pStyle w:val="Heading5"/></w:pPr>'
$inner = $inner + '<w:bookmarkStart w:id="4" w:name="section"/><w:bookmarkEnd w:id="4"/>'
Set-ParagraphXml 11 $inner

# --- Merge the two runs of the "Markdown code" paragraph into one -------
$quoteOpen = [char]0x2018
$quoteClose = [char]0x2019
$inner = '<w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr>'
$inner = $inner + '<w:r><w:t>The Markdown code ' + $quoteOpen + '#####' + $quoteClose + ' is originally reserved for header 5. However, we will use it to insert a pagebreak in a .docx document.</w:t></w:r>'
Set-ParagraphXml 12 $inner
